$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 11 ("(End Session 2)") to make room
# for a new "Open Mike" slot for session 2, shifting everything below down
# by one row (old row 11 -> new row 12, ..., old row 24 -> new row 25).
$ws.Rows("11:11").Insert()

# Copy the formatting of the row above (11:00 "Slides & Videos..." row) into
# the newly inserted row so borders/number-format/fonts match the rest of
# the schedule table.
$ws.Range("B10:D10").Copy()
$ws.Range("B11:D11").PasteSpecial(-4122)

# Fill in the new row's content: 11:30 AM, "Open Mike".
$ws.Range("B11").Value = 0.47916666666666669
$ws.Range("C11").Value = "Open Mike"

# Match the author's final cursor position.
$ws.Range("G30").Select() | Out-Null
